$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "20.193.27"
Set-TextValue "E2" "  +1.26%  "
Set-TextValue "D3" "1.440.13"
Set-TextValue "E3" "  +2.37%  "
Set-TextValue "E4" "  +0.62%  "
Set-TextValue "D5" "0.9185"
Set-TextValue "E5" "  -8.27%  "
Set-TextValue "D6" "277.03"
Set-TextValue "E6" "  +2.20%  "
Set-TextValue "D7" "0.3666"
Set-TextValue "E7" "  -0.84%  "
Set-TextValue "D8" "0.3130"
Set-TextValue "E8" "  +2.39%  "
Set-TextValue "D9" "38.87"
Set-TextValue "E9" "  -0.24%  "
Set-TextValue "E10" "  +4.87%  "
Set-TextValue "D11" "0.06532"
Set-TextValue "E11" "  +0.49%  "
Set-TextValue "E12" "  -0.03%  "
Set-TextValue "D13" "5.408"
Set-TextValue "E13" "  +2.01%  "
Set-TextValue "E14" "  +5.46%  "
Set-TextValue "D15" "6.070"
Set-TextValue "E15" "  -0.38%  "
Set-TextValue "D16" "1.445.61"
Set-TextValue "E16" "  +2.75%  "
Set-TextValue "D17" "0.00001015"
Set-TextValue "E17" "  +1.08%  "
Set-TextValue "D18" "0.9348"
Set-TextValue "E18" "  -6.63%  "
Set-TextValue "D19" "0.05641"
Set-TextValue "E19" "  -0.91%  "
Set-TextValue "D20" "67.98"
Set-TextValue "E20" "  -6.03%  "
Set-TextValue "E21" "  -1.65%  "
Set-TextValue "D22" "14.49"
Set-TextValue "E22" "  +1.98%  "
Set-TextValue "D23" "10.86"
Set-TextValue "E23" "  +1.50%  "
Set-TextValue "D24" "2.274"
Set-TextValue "E24" "  +0.15%  "
Set-TextValue "D25" "20.222.05"
Set-TextValue "E25" "  +1.33%  "
Set-TextValue "D26" "2.191"
Set-TextValue "E26" "  -0.21%  "
Set-TextValue "D27" "136.43"
Set-TextValue "E27" "  +0.19%  "
Set-TextValue "D28" "16.98"
Set-TextValue "E28" "  +2.34%  "
Set-TextValue "D29" "1.597.43"
Set-TextValue "E29" "  +2.06%  "
Set-TextValue "D30" "110.89"
Set-TextValue "E30" "  +2.82%  "
Set-TextValue "D31" "3.766"
Set-TextValue "E31" "  -2.16%  "
Set-TextValue "D32" "0.8119"
Set-TextValue "E32" "  +1.17%  "
Set-TextValue "D33" "4.834"
Set-TextValue "E33" "  -7.30%  "
Set-TextValue "D34" "0.07698"
Set-TextValue "E34" "  +0.34%  "
Set-TextValue "D35" "0.06027"
Set-TextValue "E35" "  +4.12%  "
Set-TextValue "D36" "1.484"
Set-TextValue "E36" "  +13.57%  "
Set-TextValue "D37" "4.725"
Set-TextValue "E37" "  +0.01%  "
Set-TextValue "D38" "1.134"
Set-TextValue "E38" "  +6.54%  "
Set-TextValue "B39" "VeChain"
Set-TextValue "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.01994"
Set-TextValue "E39" "  -1.07%  "
Set-TextValue "B40" "Aptos"
Set-TextValue "C40" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D40" "10.22"
Set-TextValue "E40" "  +1.36%  "
Set-TextValue "D41" "0.9351"
Set-TextValue "E41" "  -6.57%  "
Set-TextValue "D42" "0.1835"
Set-TextValue "E42" "  -4.89%  "
Set-TextValue "D43" "7.248"
Set-TextValue "E43" "  -13.13%  "
Set-TextValue "D44" "0.5254"
Set-TextValue "E44" "  +0.03%  "
Set-TextValue "E45" "  +0.58%  "
Set-TextValue "D46" "12.00"
Set-TextValue "E46" "  -0.37%  "
Set-TextValue "D47" "119.49"
Set-TextValue "E47" "  +9.06%  "
Set-TextValue "D48" "0.5156"
Set-TextValue "E48" "  +1.89%  "
Set-TextValue "D49" "1.772"
Set-TextValue "E49" "  +0.19%  "
Set-TextValue "D50" "0.06330"
Set-TextValue "E50" "  +2.97%  "
Set-TextValue "D51" "0.9940"
Set-TextValue "E51" "  -0.74%  "
